$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the (1-based) Paragraphs index of the paragraph whose text
# contains $anchorText. Throws if not found so failures are loud, not silent.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($anchorText) {
    $count = $d.Paragraphs.Count
    $i = 1
    while ($i -le $count) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like "*$anchorText*") {
            return $i
        }
        $i = $i + 1
    }
    throw "paragraph not found for anchor: $anchorText"
}

# ---------------------------------------------------------------------------
# Strike the whole paragraph (its run text AND its paragraph mark), matching
# the diff pattern where trailing "Criterios de aceite" / helper paragraphs
# are struck through completely.
# ---------------------------------------------------------------------------
function Strike-WholeParagraph($anchorText) {
    $idx = Find-ParagraphIndex $anchorText
    $p = $d.Paragraphs($idx)
    $p.Range.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# Strike everything in the paragraph AFTER the literal "Prompt: " prefix,
# leaving "Prompt: " itself unformatted - matching the diff pattern where the
# leading run is split into an unstruck "Prompt: " run followed by one (or
# more, for paragraphs that already contained VerbatimChar-styled runs) struck
# runs covering the remainder of the paragraph text. The paragraph mark itself
# is intentionally left untouched.
# ---------------------------------------------------------------------------
function Strike-AfterPrompt($anchorText) {
    $idx = Find-ParagraphIndex $anchorText
    $p = $d.Paragraphs($idx)

    $searchRng = $p.Range.Duplicate
    $found = $searchRng.Find.Execute("Prompt: ")
    if (-not $found) {
        throw "'Prompt: ' prefix not found in paragraph: $anchorText"
    }

    $afterPromptStart = $searchRng.End
    $paragraphEnd = $p.Range.End

    $rest = $d.Range($afterPromptStart, $paragraphEnd)
    $rest.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# Prompt 2 - "No pacote core, defina..." block
# ---------------------------------------------------------------------------
Strike-AfterPrompt "No pacote core, defina:"
Strike-WholeParagraph "Compila - Há testes unitários vazios preparados para essas funções"

# ---------------------------------------------------------------------------
# Prompt 3 - "Adicione em src/core/normalizers.ts..." block
# ---------------------------------------------------------------------------
Strike-AfterPrompt "Adicione em src/core/normalizers.ts:"
Strike-WholeParagraph "Crie testes vitest para cada função com casos comuns e edge cases."
Strike-WholeParagraph "vitest ok - funções cobrem casos básicos e não quebram inputs estranhos"

# ---------------------------------------------------------------------------
# Prompt 4 - "Crie src/core/themeLoader.ts com..." block
# ---------------------------------------------------------------------------
Strike-AfterPrompt "Crie src/core/themeLoader.ts com:"
Strike-WholeParagraph "Crie testes simulando temas simples."
Strike-WholeParagraph "Carrega tema quando existir - Fallback funciona sem tema"

# ---------------------------------------------------------------------------
# Prompt 5 - "Implemente em src/core/matchers..." block
# ---------------------------------------------------------------------------
Strike-AfterPrompt "Implemente em src/core/matchers:"
Strike-WholeParagraph "Inclua testes: entradas em px, rem, %, números, e shorthands."
Strike-WholeParagraph "Casos básicos passam - Saída consistente"

Write-Output "applied strikethrough to Prompt 2/3/4/5 blocks"
